$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "ok"
$ws.Range("B2").Value = "ok"
$ws.Range("B3").Value = "ok"
$ws.Range("B4").Value = "ok"
$ws.Range("B5").Value = "perlu disesuaikan dgn kebutuhan"
$ws.Range("B6").Value = "perlu disesuaikan dgn kebutuhan (untuk hapus data)"

$ws.Range("M10").Select()
